$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update sheet1 ("o_10"): add column E, update row 2 values ---

$promptText = @"
 Given is the adjacency matrix for a weighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the least cost path from node A to node I?
   A B C D E F G H I
 A 0 1 0 3 0 0 0 0 0
 B 1 0 2 0 2 0 0 0 0
 C 0 2 0 0 0 2 0 0 0
 D 3 0 0 0 1 0 2 0 0
 E 0 2 0 1 0 3 0 1 0
 F 0 0 2 0 3 0 0 0 1
 G 0 0 0 2 0 0 0 2 0
 H 0 0 0 0 1 0 2 0 1
 I 0 0 0 0 0 1 0 1 0

Solution: A -> B -> E -> H -> I
        

Example 2: what is the least cost path from node A to node I?
   A B C D E F G H I
 A 0 4 0 2 0 0 0 0 0
 B 4 0 4 0 3 0 0 0 0
 C 0 4 0 0 0 3 0 0 0
 D 2 0 0 0 0 0 4 0 0
 E 0 3 0 0 0 3 0 0 0
 F 0 0 3 0 3 0 0 0 1
 G 0 0 0 4 0 0 0 3 0
 H 0 0 0 0 0 0 3 0 4
 I 0 0 0 0 0 1 0 4 0

Solution: A -> B -> E -> F -> I
        

Example 3: what is the least cost path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 4 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 5 0 0 4 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 3 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 2 0 0 4 0 0 0 0 0 0 0
 F 0 0 0 0 2 0 4 0 0 0 0 0 0 0 0 0
 G 0 0 4 0 0 4 0 5 0 0 5 0 0 0 0 0
 H 0 0 0 3 0 0 5 0 0 0 0 5 0 0 0 0
 I 0 0 0 0 4 0 0 0 0 2 0 0 4 0 0 0
 J 0 0 0 0 0 0 0 0 2 0 0 0 0 4 0 0
 K 0 0 0 0 0 0 5 0 0 0 0 2 0 0 1 0
 L 0 0 0 0 0 0 0 5 0 0 2 0 0 0 0 4
 M 0 0 0 0 0 0 0 0 4 0 0 0 0 5 0 0
 N 0 0 0 0 0 0 0 0 0 4 0 0 5 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 5 0 3
 P 0 0 0 0 0 0 0 0 0 0 0 4 0 0 3 0

Solution: A -> E -> F -> G -> K -> O -> P
        
 Given these examples, answer the following quesiton.

what is the least cost path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 1 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 2 0 0 2 0 0 0 0 0 0 0 0 0 0
 C 0 2 0 1 0 0 5 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 5 0 0 0 0 0 0 0 0
 E 3 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0
 F 0 2 0 0 0 0 4 0 0 1 0 0 0 0 0 0
 G 0 0 5 0 0 4 0 3 0 0 5 0 0 0 0 0
 H 0 0 0 5 0 0 3 0 0 0 0 3 0 0 0 0
 I 0 0 0 0 2 0 0 0 0 4 0 0 1 0 0 0
 J 0 0 0 0 0 1 0 0 4 0 0 0 0 3 0 0
 K 0 0 0 0 0 0 5 0 0 0 0 0 0 0 4 0
 L 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 5
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 5 0 0
 N 0 0 0 0 0 0 0 0 0 3 0 0 5 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 4 0 0 5 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 5 0 0 1 0
    
"@

$ws1.Cells.Item(1,5).Value = "evaluator_partial_correctness"
$ws1.Cells.Item(2,1).Value = $promptText
$ws1.Cells.Item(2,2).Value = "A -> B -> F -> J -> N -> O -> P"
$ws1.Cells.Item(2,3).Value = "The least cost path from node A to node P is A -> B -> C -> G -> O -> P."
$ws1.Cells.Item(2,4).Value = "Wrong"
$ws1.Cells.Item(2,5).Value = "Output: 3/6"

# Copy header style (from D1, which already has the bold/border/center style) onto E1
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)

# Undo the implicit row-autofit-on-assignment so row 2 keeps the default height
$ws1.Rows.Item(2).AutoFit()

# --- Add new sheets "o_20" and "o_20_jumbled" after "o_10" ---

$headers = @("prompt", "solution", "llm_response", "evaluator_response", "evaluator_partial_correctness")

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"
for ($c = 1; $c -le 5; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c-1]
}
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"
for ($c = 1; $c -le 5; $c++) {
    $ws3.Cells.Item(1, $c).Value = $headers[$c-1]
}
$ws1.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)

$ws1.Select() | Out-Null
